# Append two new match rows (73 and 74) to the Thai League 1 2023-2024 sheet,
# matching the format already used for the rest of the data (bold/bordered
# index column, date-time formatted match-date column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row down onto the
# two new rows before filling in values, so the new cells keep the same
# look (A column bold/border style, E column date-time number format).
$ws.Range("A72:V72").Copy()
$ws.Range("A73:V74").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{
        Row = 73
        Values = @(
            72, "thailand", "thai-league-1", "2023-2024", 45235.47916666666,
            "Nakhon Pathom", 2, "Sukhothai", 2,
            1.81, "04/11/2023 17:13",
            2.27, "05/11/2023 11:25",
            3.93, "04/11/2023 17:13",
            3.46, "05/11/2023 11:22",
            3.74, "04/11/2023 17:13",
            3.18, "05/11/2023 11:25",
            "https://www.betexplorer.com/football/thailand/thai-league-1/nakhon-pathom-sukhothai/CMI1gBVn/"
        )
    },
    @{
        Row = 74
        Values = @(
            73, "thailand", "thai-league-1", "2023-2024", 45235.54166666666,
            "Ratchaburi", 1, "Chonburi", 2,
            1.95, "31/10/2023 10:42",
            2.11, "05/11/2023 12:56",
            3.63, "31/10/2023 10:42",
            3.5, "05/11/2023 12:56",
            3.72, "31/10/2023 10:42",
            3.5, "05/11/2023 12:56",
            "https://www.betexplorer.com/football/thailand/thai-league-1/ratchaburi-chonburi/6qdzwjo6/"
        )
    }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $col = 1
    foreach ($val in $entry.Values) {
        $ws.Cells.Item($r, $col).Value = $val
        $col++
    }
}
